$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new test-log entry as row 5 (date 2017-08-06 = serial 42953).
# Copy formatting from existing rows so the same cell styles (date format,
# wrap-text "Zu Verbessern" column, plain link columns) are reused instead
# of new styles being synthesized.
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D4").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Range("E3").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("I5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A5").Value = 42953
$ws.Range("D5").Value = "Sichtfeld: mit TriggerCollider ODER mit Raycasting/Angle?"
$ws.Range("F5").Value = "https://gamedev.stackexchange.com/questions/104773/vision-cone-for-enemy-ai-in-unity-2d"
$ws.Range("G5").Value = "http://answers.unity3d.com/questions/414479/2d-enemy-field-of-vision-script.html"
$ws.Range("H5").Value = "https://www.youtube.com/watch?v=rQG9aUWarwE"
$ws.Range("I5").Value = "https://www.youtube.com/watch?v=mBGUY7EUxXQ"

$ws.Rows.Item(5).RowHeight = 30

$ws.Range("F11").Select()
